# doc: tech implement path
# Add two new worksheets (Sheet3, Sheet4) describing the step-by-step /
# parallel-task account-balance-transfer walkthrough, and update the
# workbook's view state so that the new last sheet (Sheet4) is the active one.

$wb = $excel.ActiveWorkbook

# Best-effort: reflect the saved window size from the target workbook view.
$mainWindow = $wb.Windows.Item(1)
$mainWindow.Width = 27660
$mainWindow.Height = 13160

# ---------------------------------------------------------------------------
# 1. Create Sheet3 and Sheet4 at the end of the workbook (after current last
#    sheet), in order, so they are named "Sheet3" / "Sheet4" automatically.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3 = $wb.Worksheets.Add($null, $lastSheet)

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4 = $wb.Worksheets.Add($null, $lastSheet)

# ---------------------------------------------------------------------------
# 2. Fill in Sheet3 ("serial" walkthrough) data, row by row / left to right
#    so that new shared strings are created in a stable, deterministic order.
# ---------------------------------------------------------------------------
$sheet3Data = @(
    @("步骤", "A", "B", "C", "备注"),
    @(1, "< 2", "/", "< 0", "从账本中读取A的账户余额和C的账户余额"),
    @(2, "2-0.5", "/", "0+0.5", "计算收支平衡：A账户余额减去0.5，C账户余额增加0.5"),
    @(3, "> 1.5", "/", "> 0.5", "将A和C的账户余额更新回账本"),
    @(4, "/", "< 1", "< 0.5", "从账本中读取B的账户余额和C的账户余额"),
    @(5, "/", "1-0.5", "0.5+0.5", "计算收支平衡：B账户余额减去0.5，C账户余额增加 0.5"),
    @(6, "/", "> 0.5", "> 1", "将B和C的账户余额更新回账本")
)

for ($r = 0; $r -lt $sheet3Data.Length; $r++) {
    $row = $sheet3Data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sheet3.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------------
# 3. Fill in Sheet4 ("parallel tasks" walkthrough) data the same way.
# ---------------------------------------------------------------------------
$sheet4Data = @(
    @("步骤", "A", "B", "C", "备注"),
    @(1, "< 2", "/", "< 0", "并行任务1：从账本中读取A和C的账户余额"),
    @(2, "2-0.5", "/", "0+0.5", "并行任务1：计算收支平衡，A账户余额减去0.5，C账户余额增加0.5"),
    @("*3", "/", "< 1", "< 0", "并行任务2：从账本中读取B和C的账户余额（注意，还是0！）"),
    @(4, "> 1.5", "/", "> 0.5", "并行任务1：将A和C的账户余额更新回账本"),
    @("*5", "/", "1-0.5", "0+0.5", "并行任务2：计算收支平衡，B账户余额减去0.5，C账户余额增加0.5"),
    @("*6", "/", "> 0.5", "> 0.5", "并行任务2：将B和C的账户余额更新回账本")
)

for ($r = 0; $r -lt $sheet4Data.Length; $r++) {
    $row = $sheet4Data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $sheet4.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# ---------------------------------------------------------------------------
# 4. Widen column E (the "备注"/remarks column) on both new sheets so the
#    long descriptive text is readable.
# ---------------------------------------------------------------------------
$sheet3.Columns.Item(5).ColumnWidth = 52
$sheet4.Columns.Item(5).ColumnWidth = 70.4

# ---------------------------------------------------------------------------
# 5. Recreate the leftover outline-level bookkeeping (outlineLevelRow /
#    outlineLevelCol) that is present on both new sheets.
# ---------------------------------------------------------------------------
$sheet3.Rows.Item(7).OutlineLevel = 6
$sheet3.Columns.Item(5).OutlineLevel = 4
$sheet4.Rows.Item(7).OutlineLevel = 6
$sheet4.Columns.Item(5).OutlineLevel = 4

# ---------------------------------------------------------------------------
# 6. Restore per-sheet selections.
# ---------------------------------------------------------------------------
$sheet3.Activate()
[void]$sheet3.Range("E4").Select()

$sheet4.Activate()
[void]$sheet4.Range("E2:E7").Select()

# ---------------------------------------------------------------------------
# 7. Leave Sheet4 (the new last sheet) as the active / selected tab, matching
#    the saved workbook view (activeTab points at the 4th sheet).
# ---------------------------------------------------------------------------
$sheet4.Activate()
